# Update the "Correspond Handoff Datetime" (col E) and
# "Correspond Handback DateTime" (col H) values for the
# "8e14a4a7-4766-4ae2-a79c-de85c5dc3078..." rows on the zh-cn and de-de
# sheets, reflecting regenerated report timestamps.

$wb = $excel.ActiveWorkbook

# zh-cn sheet: row 3 and its duplicate row 5 share the same values
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "2016-03-22 06:23:40"
$wsZhCn.Range("E5").Value = "2016-03-22 06:23:40"
$wsZhCn.Range("H3").Value = "2016-03-22 06:24:04"
$wsZhCn.Range("H5").Value = "2016-03-22 06:24:04"

# de-de sheet: row 3 and its duplicate row 5 share the same values
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "2016-03-22 06:23:45"
$wsDeDe.Range("E5").Value = "2016-03-22 06:23:45"
$wsDeDe.Range("H3").Value = "2016-03-22 06:24:10"
$wsDeDe.Range("H5").Value = "2016-03-22 06:24:10"
